$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Dlk1", "Notch1", "ECs", 1, 0.3333333333333333, 0.028147, 0.084441, 0.0007347845853173872, 0.000734784585317387, 3, 1, 68.00339, 204.01017, 0.6265962299909886, 0.6265962299909885, 1.91409141833, 17.22682276497, 0.0004604132510153667, 0.0004604132510153666)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i+1).Value = $row2[$i] }

$row3 = @("ECs", "Dlk1", "Notch1", "FAPs", 1, 0.3333333333333333, 0.028147, 0.084441, 0.0007347845853173872, 0.000734784585317387, 3, 1, 8.99153, 26.97459, 0.08284967558015671, 0.08284967558015671, 0.25308459491, 2.27776135419, 0.0000608766645148455, 0.0000608766645148455)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i+1).Value = $row3[$i] }

$row4 = @("ECs", "Dlk1", "Notch1", "M2", 1, 0.3333333333333333, 0.028147, 0.084441, 0.0007347845853173872, 0.000734784585317387, 3, 1, 10.76843933333333, 32.305318, 0.09922245772090688, 0.09922245772090688, 0.3030992619153334, 2.727893357238, 0.00007290713245062853, 0.00007290713245062853)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i+1).Value = $row4[$i] }

$row5 = @("ECs", "Dlk1", "Notch1", "sCs", 1, 0.3333333333333333, 0.028147, 0.084441, 0.0007347845853173872, 0.000734784585317387, 3, 1, 20.764887, 62.294661, 0.1913316367079478, 0.1913316367079478, 0.584469274389, 5.260223469501, 0.0001405875373365464, 0.0001405875373365464)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i+1).Value = $row5[$i] }

$row6 = @("FAPs", "Dlk1", "Notch1", "ECs", 3, 1, 36.97491766666666, 110.924753, 0.9652396186039789, 0.9652396186039789, 3, 1, 68.00339, 204.01017, 0.6265962299909886, 0.6265962299909885, 2514.419746304223, 22629.77771673801, 0.6048155060551929, 0.6048155060551927)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i+1).Value = $row6[$i] }

$row7 = @("FAPs", "Dlk1", "Notch1", "FAPs", 3, 1, 36.97491766666666, 110.924753, 0.9652396186039789, 0.9652396186039789, 3, 1, 8.99153, 26.97459, 0.08284967558015671, 0.08284967558015671, 332.4610814473633, 2992.14973302627, 0.07996978925845384, 0.07996978925845384)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i+1).Value = $row7[$i] }

$row8 = @("FAPs", "Dlk1", "Notch1", "M2", 3, 1, 36.97491766666666, 110.924753, 0.9652396186039789, 0.9652396186039789, 3, 1, 10.76843933333333, 32.305318, 0.09922245772090688, 0.09922245772090688, 398.1621577484948, 3583.459419736454, 0.09577344724747758, 0.09577344724747758)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i+1).Value = $row8[$i] }

$row9 = @("FAPs", "Dlk1", "Notch1", "sCs", 3, 1, 36.97491766666666, 110.924753, 0.9652396186039789, 0.9652396186039789, 3, 1, 20.764887, 62.294661, 0.1913316367079478, 0.1913316367079478, 767.7799871826369, 6910.019884643732, 0.1846808760428546, 0.1846808760428546)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i+1).Value = $row9[$i] }

$row10 = @("sCs", "Dlk1", "Notch1", "ECs", 3, 1, 1.303400333333333, 3.910201, 0.03402559681070371, 0.03402559681070371, 3, 1, 68.00339, 204.01017, 0.6265962299909886, 0.6265962299909885, 88.63564119379666, 797.72077074417, 0.02132031068478035, 0.02132031068478035)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i+1).Value = $row10[$i] }

$row11 = @("sCs", "Dlk1", "Notch1", "FAPs", 3, 1, 1.303400333333333, 3.910201, 0.03402559681070371, 0.03402559681070371, 3, 1, 8.99153, 26.97459, 0.08284967558015671, 0.08284967558015671, 11.71956319917667, 105.47606879259, 0.002819009657188017, 0.002819009657188017)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, $i+1).Value = $row11[$i] }

$row12 = @("sCs", "Dlk1", "Notch1", "M2", 3, 1, 1.303400333333333, 3.910201, 0.03402559681070371, 0.03402559681070371, 3, 1, 10.76843933333333, 32.305318, 0.09922245772090688, 0.09922245772090688, 14.03558741654644, 126.320286748918, 0.003376103340978673, 0.003376103340978673)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, $i+1).Value = $row12[$i] }

$row13 = @("sCs", "Dlk1", "Notch1", "sCs", 3, 1, 1.303400333333333, 3.910201, 0.03402559681070371, 0.03402559681070371, 3, 1, 20.764887, 62.294661, 0.1913316367079478, 0.1913316367079478, 27.064960637429, 243.584645736861, 0.006510173127756669, 0.006510173127756669)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, $i+1).Value = $row13[$i] }
